# Apply layout/font tweaks to the title slide (slide 1) of the presentation:
#  - Title placeholder: reposition/resize and shrink font from 34pt to 32pt
#  - Subtitle placeholder: move down slightly
#  - Author block shape: move down slightly
#
# Note: PowerPoint Shape.Left/Top/Width/Height are expressed in points;
# a few target values below use the nearest point value whose internal
# (single-precision) representation converts back to the exact target
# EMU offset from the reference OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1: Title ("Rectangle 2" / ctrTitle placeholder)
$title = $s.Shapes.Item(1)
$title.Left = 18
$title.Top = 34.5
$title.Width = 684
$title.Height = 118.12496185302734
$title.TextFrame.TextRange.Font.Size = 32

# Shape 2: Subtitle ("Rectangle 3" / subTitle placeholder)
$subtitle = $s.Shapes.Item(2)
$subtitle.Left = 57
$subtitle.Top = 152.62496948242188
$subtitle.Width = 606
$subtitle.Height = 54

# Shape 3: Author block ("Rectangle 4")
$authors = $s.Shapes.Item(3)
$authors.Left = 108
$authors.Top = 210.65614318847656
$authors.Width = 570
$authors.Height = 112.3125991821289
